$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 3 for columns D, M, N, O, P, S
$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $r2 = $ws.Range($col + "2")
    $r3 = $ws.Range($col + "3")
    $v2 = $r2.Value2
    $v3 = $r3.Value2
    $r2.Value2 = $v3
    $r3.Value2 = $v2
}
